$d = $word.ActiveDocument

# 1. Remove the old front-matter paragraphs (Title, Overview heading,
#    the two body paragraphs, the "Test Code (Python)" heading, and the
#    "import re" line) - six paragraphs total at the top of the document.
$deleteStart = $d.Paragraphs.Item(1).Range.Start
$deleteEnd = $d.Paragraphs.Item(6).Range.End
$d.Range($deleteStart, $deleteEnd).Delete()

# 2. Append the new code block (from the Test Lead) as a series of
#    Title-styled paragraphs at the end of the document.
$newParagraphs = @(
    'Test Code Document – Password Validation Logic',
    'def validate_password(password):',
    '  if len(password) < 8:',
    '    return False, "Password must be at least 8 characters long"',
    '  if not re.search(r"\d", password):',
    '    return False, "Password must contain at least one digit"',
    '  if not re.search(r"[!@#$%^&*(),.?\":{}|<>`]", password):',
    '    return False, "Password must contain at least one special character"',
    '  return True, "Password is valid"',
    'def test_passwords():',
    '  print(valid, msg)',
    '  print(valid, msg)',
    '  print(valid, msg)',
    '  print(valid, msg)',
    '  print(valid, msg)',
    '  print(valid, msg)',
    '  print(valid, msg)',
    'if __name__ == "__main__":'
)

# RGB(0, 176, 80) == hex 00B050, Word stores BGR-packed OLE_COLOR as r + g*256 + b*65536
$greenColor = 0 + (176 * 256) + (80 * 65536)

for ($i = 0; $i -lt $newParagraphs.Count; $i++) {
    $tail = $d.Paragraphs.Item($d.Paragraphs.Count).Range
    $tail.InsertParagraphAfter()
    $newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
    $newPara.Style = "Title"
    $newPara.Range.Text = $newParagraphs[$i]
    $textRange = $d.Range($newPara.Range.Start, $newPara.Range.End - 1)
    if ($i -eq 0) {
        # First new paragraph (repeated title) keeps default title coloring;
        # toggling Bold on/off mirrors the leftover empty run-properties left
        # behind by the author clearing direct formatting.
        $textRange.Font.Bold = 1
        $textRange.Font.Bold = 0
    } else {
        $textRange.Font.Color = $greenColor
    }
}

Write-Host "Final paragraph count:" $d.Paragraphs.Count
